$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.727399999999996
$ws.Range("A8").Value = -22.50070000000002
$ws.Range("A10").Value = -21.85319999999998
$ws.Range("A12").Value = -21.48859999999998
$ws.Range("B12").Value = 5.5612
$ws.Range("D12").Value = -5.911900000000003
$ws.Range("D13").Value = -8.807899999999989
$ws.Range("B15").Value = 5.680099999999994
$ws.Range("B17").Value = 4.532300000000001
$ws.Range("A18").Value = -22.33970000000002
$ws.Range("D21").Value = -7.53439999999999
$ws.Range("D25").Value = -7.346099999999996
$ws.Range("B26").Value = 4.439000000000001
$ws.Range("B27").Value = 6.512100000000003
$ws.Range("B28").Value = 6.284199999999998
$ws.Range("D32").Value = -7.664699999999995
$ws.Range("D36").Value = -7.287499999999998
$ws.Range("A37").Value = -21.90679999999999
$ws.Range("B37").Value = 6.462100000000003
$ws.Range("D38").Value = -8.099800000000004
$ws.Range("D41").Value = -8.437499999999996
$ws.Range("B47").Value = 6.747300000000004
$ws.Range("D52").Value = -7.799000000000003
$ws.Range("A55").Value = -22.07230000000001
$ws.Range("D59").Value = -8.681599999999994
$ws.Range("B65").Value = 5.376700000000003
$ws.Range("D67").Value = -7.050499999999996
$ws.Range("A68").Value = -21.4749
$ws.Range("B73").Value = 8.769599999999999
$ws.Range("A77").Value = -19.81369999999999
$ws.Range("A78").Value = -19.77669999999999
$ws.Range("A81").Value = -22.20210000000001
$ws.Range("A82").Value = -21.9004
$ws.Range("B84").Value = 5.907200000000001
$ws.Range("D84").Value = -7.7149
$ws.Range("B85").Value = 5.978300000000001
$ws.Range("D88").Value = -8.0692
$ws.Range("D89").Value = -8.069699999999997
$ws.Range("B93").Value = 5.457600000000002
$ws.Range("B95").Value = 5.234900000000001
$ws.Range("D95").Value = -7.559699999999997
$ws.Range("B98").Value = 5.816000000000005
$ws.Range("B99").Value = 6.0736
$ws.Range("B101").Value = 6.456499999999998
$ws.Range("D105").Value = -8.301200000000005
